# Automatische test-sync: 2025-06-29 15:11:50
# Append the new "Testmail #14" log entry to the "Logs" sheet (row 29),
# extend the conditional-formatting ranges to include the new row, and
# bump the "Bestelling / Levering" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 29

$logs.Cells.Item($newRow, 1).Value = "Kun je deze bestelling vandaag verwerken?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #14: Kun je deze bestelling vandaag verwerken?"
$logs.Cells.Item($newRow, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,
Hartelijk dank voor uw e-mail. Om de bestelling vandaag te kunnen verwerken, hebben we een ordernummer of klantgegevens nodig. Zou u ons deze gegevens kunnen verstrekken, zodat we uw verzoek verder kunnen onderzoeken?
Met vriendelijke groet,
[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-06-29 15:11:04"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"

# Writing the multi-line reply text into column E causes the row to pick
# up an explicit (wrong) auto-computed height; AutoFit puts the row back
# into "no explicit height" state, matching the other data rows.
$logs.Rows.Item($newRow).AutoFit()

# Extend the existing conditional formatting blocks so they keep covering
# every data row (they previously stopped at row 28).
$columns = @("D", "G", "H", "I")
foreach ($col in $columns) {
    $oldRange = $logs.Range("${col}2:${col}28")
    $newRange = $logs.Range("${col}2:${col}${newRow}")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Bestelling / Levering" (7 -> 8).
$dashboard.Cells.Item(2, 2).Value = 8
